$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 393, shifting existing rows 393:410 down to 394:411
$ws.Rows.Item(393).EntireRow.Insert()

# Copy the constant / unchanged columns from the row that was pushed down (now row 394)
$ws.Range("A393").Value = $ws.Range("A394").Value2
$ws.Range("B393").Value = $ws.Range("B394").Value2
$ws.Range("C393").Value = $ws.Range("C394").Value2
$ws.Range("E393").Value = $ws.Range("E394").Value2
$ws.Range("F393").Value = $ws.Range("F394").Value2
$ws.Range("G393").Value = $ws.Range("G394").Value2
$ws.Range("H393").Value = $ws.Range("H394").Value2
$ws.Range("I393").Value = $ws.Range("I394").Value2
$ws.Range("N393").Value = $ws.Range("N394").Value2
$ws.Range("O393").Value = $ws.Range("O394").Value2
$ws.Range("Q393").Value = $ws.Range("Q394").Value2
$ws.Range("R393").Value = $ws.Range("R394").Value2

# Set the new data for this inserted record
$ws.Range("D393").Value = 44931
$ws.Range("J393").Value = 500
$ws.Range("K393").Value = 1000
$ws.Range("L393").Value = 1000
$ws.Range("M393").Value = 1000
$ws.Range("P393").Value = 200
